$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 115.4136146666667
$ws.Range("H2").Value = 346.240844
$ws.Range("I2").Value = 0.2619217538490851
$ws.Range("J2").Value = 0.2619217538490851
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1837056666666667
$ws.Range("N2").Value = 0.5511170000000001
$ws.Range("O2").Value = 0.7269991860920679
$ws.Range("P2").Value = 0.7269991860920678
$ws.Range("Q2").Value = 21.20213502474978
$ws.Range("R2").Value = 190.819215222748
$ws.Range("S2").Value = 0.1904169018680918
$ws.Range("T2").Value = 0.1904169018680918

# Row 3
$ws.Range("G3").Value = 115.4136146666667
$ws.Range("H3").Value = 346.240844
$ws.Range("I3").Value = 0.2619217538490851
$ws.Range("J3").Value = 0.2619217538490851
$ws.Range("O3").Value = 0.2534828531892131
$ws.Range("P3").Value = 0.2534828531892131
$ws.Range("Q3").Value = 7.39254978903911
$ws.Range("R3").Value = 66.53294810135199
$ws.Range("S3").Value = 0.06639267347798884
$ws.Range("T3").Value = 0.06639267347798886

# Row 4
$ws.Range("G4").Value = 115.4136146666667
$ws.Range("H4").Value = 346.240844
$ws.Range("I4").Value = 0.2619217538490851
$ws.Range("J4").Value = 0.2619217538490851
$ws.Range("O4").Value = 0.01951796071871896
$ws.Range("P4").Value = 0.01951796071871896
$ws.Range("Q4").Value = 0.5692199475359999
$ws.Range("R4").Value = 5.122979527824
$ws.Range("S4").Value = 0.005112178503004418
$ws.Range("T4").Value = 0.005112178503004419

# Row 5
$ws.Range("I5").Value = 0.6414314537852458
$ws.Range("J5").Value = 0.6414314537852458
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1837056666666667
$ws.Range("N5").Value = 0.5511170000000001
$ws.Range("O5").Value = 0.7269991860920679
$ws.Range("P5").Value = 0.7269991860920678
$ws.Range("Q5").Value = 51.92282081354823
$ws.Range("R5").Value = 467.3053873219341
$ws.Range("S5").Value = 0.4663201448357256
$ws.Range("T5").Value = 0.4663201448357255

# Row 6
$ws.Range("I6").Value = 0.6414314537852458
$ws.Range("J6").Value = 0.6414314537852458
$ws.Range("O6").Value = 0.2534828531892131
$ws.Range("P6").Value = 0.2534828531892131
$ws.Range("S6").Value = 0.162591875030789
$ws.Range("T6").Value = 0.162591875030789

# Row 7
$ws.Range("I7").Value = 0.6414314537852458
$ws.Range("J7").Value = 0.6414314537852458
$ws.Range("O7").Value = 0.01951796071871896
$ws.Range("P7").Value = 0.01951796071871896
$ws.Range("S7").Value = 0.01251943391873122
$ws.Range("T7").Value = 0.01251943391873122

# Row 8
$ws.Range("I8").Value = 0.09664679236566912
$ws.Range("J8").Value = 0.09664679236566913
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.1837056666666667
$ws.Range("N8").Value = 0.5511170000000001
$ws.Range("O8").Value = 0.7269991860920679
$ws.Range("P8").Value = 0.7269991860920678
$ws.Range("Q8").Value = 7.823398825538333
$ws.Range("R8").Value = 70.410589429845
$ws.Range("S8").Value = 0.07026213938825053
$ws.Range("T8").Value = 0.07026213938825053

# Row 9
$ws.Range("I9").Value = 0.09664679236566912
$ws.Range("J9").Value = 0.09664679236566913
$ws.Range("O9").Value = 0.2534828531892131
$ws.Range("P9").Value = 0.2534828531892131
$ws.Range("S9").Value = 0.02449830468043527
$ws.Range("T9").Value = 0.02449830468043527

# Row 10
$ws.Range("I10").Value = 0.09664679236566912
$ws.Range("J10").Value = 0.09664679236566913
$ws.Range("O10").Value = 0.01951796071871896
$ws.Range("P10").Value = 0.01951796071871896
$ws.Range("S10").Value = 0.001886348296983317
$ws.Range("T10").Value = 0.001886348296983317
